$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-17 03:22:01"
$wsZhCn.Range("H2").Value = "2016-03-17 03:22:41"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-17 03:22:09"
$wsDeDe.Range("H2").Value = "2016-03-17 03:22:54"
